$wb = $excel.ActiveWorkbook

# The new "Spain" sheet is a copy of "Italy" (same layout/styles), placed
# right after it, with its own market name / part number and a few cosmetic
# tweaks (column widths, row heights, active-cell selection).
$italy = $wb.Worksheets.Item("Italy")
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item("Italy (2)")
$spain.Name = "Spain"

# Spain-specific values
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2034 "

# Column widths: narrower B (market name) / D (input value) columns.
$spain.Columns.Item(2).ColumnWidth = 15.6640625
$spain.Columns.Item(4).ColumnWidth = 19.33203125

# Row 13 reverts to the sheet's standard (auto) height.
$spain.Range("13:13").EntireRow.AutoFit()

# Rows 3-5 grow to fit the wrapped "Input/Expected value" labels.
$spain.Range("3:5").RowHeight = 28.8

# Spain becomes the active tab, with B4 selected.
$spain.Activate()
$spain.Range("B4").Select()

# Italy is no longer the active tab; its last selection moves to B2.
$italy.Range("B2").Select()
$spain.Activate()
